# Fix redundant empty lines: the last data column (AQ) on the Energy,
# Entropy and Hurst sheets was a stray/redundant column. Its values are
# dropped and the previous column (AP) is updated to hold the corrected
# figures, then the now-empty AQ column is removed so the sheet's used
# range shrinks back from B1:AQ11 to B1:AP11.

$wb = $excel.ActiveWorkbook

$newApValues = @{
    "Energy" = @{
        2  = 1.480808714107192
        3  = 0.8075663965570536
        4  = 0.4025583646820428
        5  = 0.6899572908529897
        6  = 0.5377756708469189
        7  = 4.16202276582262
        8  = 0.4448370713058029
        9  = 0.6166176369715476
        10 = 0.1467987596757228
        11 = 0.2284570196298694
    }
    "Entropy" = @{
        2  = 5.81862352789444
        3  = 8.874440198440888
        4  = 8.39327524993538
        5  = 8.219518966396434
        6  = 6.815171239192523
        7  = 8.578694557752547
        8  = 8.234669674043829
        9  = 6.811470816989324
        10 = 7.450296224191622
        11 = 7.330916878114602
    }
    "Hurst" = @{
        2  = 0.6835320619646125
        3  = 0.6316767355853479
        4  = 0.5613107430178561
        5  = 0.6279557123376847
        6  = 0.5810877097758038
        7  = 0.7294717938370332
        8  = 0.5575697229131329
        9  = 0.6311833169644426
        10 = 0.648578548383285
        11 = 0.6408262317860869
    }
}

foreach ($sheetName in @("Energy", "Entropy", "Hurst")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rowValues = $newApValues[$sheetName]

    foreach ($r in 2..11) {
        $ws.Cells.Item($r, 42).Value = $rowValues[$r]
    }

    $ws.Columns("AQ").Delete()
}
